$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete the two trailing incomplete rows (old rows 23 and 24:
#    A23/B23 = 45205/east, A24/B24 = 45219/east) which only had DATE and
#    STRIP_ID filled in - no actual irrigation event data was recorded.
$ws.Rows.Item(23).EntireRow.Delete() | Out-Null
$ws.Rows.Item(23).EntireRow.Delete() | Out-Null

# 2. The WATER_DURATION..METER_GAL_USE_GAL_X_100 columns (E:J) are
#    reformatted from a mix of "h:mm" time format / unformatted General
#    to a plain "0" number format, now that duration is treated as a
#    number rather than a clock time.
$ws.Range("E1:J22").NumberFormat = "0"

# 3. Leave the cursor/selection where the author left it after their
#    cleanup pass.
$ws.Range("C29").Select() | Out-Null
